$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.949.92'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.623.83'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.58'
$ws.Range("E5").Value = '  -1.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.503'
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("E7").Value = '  +0.40%  '
$ws.Range("E8").Value = '  -2.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0618'
$ws.Range("E9").Value = '  -3.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.23'
$ws.Range("E10").Value = '  -6.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0787'
$ws.Range("E11").Value = '  -1.10%  '
$ws.Range("D12").Value = '1.848.71'
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("D13").Value = '1.619.57'
$ws.Range("E13").Value = '  -2.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.18'
$ws.Range("E14").Value = '  -2.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.524'
$ws.Range("E15").Value = '  -3.63%  '
$ws.Range("D16").Value = '25.927.81'
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.0₃0736'
$ws.Range("E17").Value = '  -3.44%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.18'
$ws.Range("E18").Value = '  -3.38%  '
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '191.31'
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.24'
$ws.Range("E21").Value = '  -2.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.57'
$ws.Range("E22").Value = '  -3.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.07'
$ws.Range("E23").Value = '  -2.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.133'
$ws.Range("E24").Value = '  +1.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.24'
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("E27").Value = '  -1.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.71'
$ws.Range("E28").Value = '  -2.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.17'
$ws.Range("E29").Value = '  -2.22%  '
$ws.Range("E30").Value = '  -1.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0481'
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.12'
$ws.Range("E32").Value = '  -4.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.11'
$ws.Range("E33").Value = '  -5.56%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.49'
$ws.Range("E34").Value = '  -2.86%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.40'
$ws.Range("E35").Value = '  -2.45%  '
$ws.Range("D36").Value = '1.124.93'
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.847'
$ws.Range("E37").Value = '  -6.22%  '
$ws.Range("E38").Value = '  -1.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.516'
$ws.Range("E39").Value = '  -4.42%  '
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.90'
$ws.Range("E41").Value = '  -1.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.766'
$ws.Range("E42").Value = '  -3.55%  '
$ws.Range("D43").Value = '1.758.90'
$ws.Range("E43").Value = '  -1.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.16'
$ws.Range("E44").Value = '  -5.66%  '
$ws.Range("E45").Value = '  -2.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '54.44'
$ws.Range("E46").Value = '  -3.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0527'
$ws.Range("E47").Value = '  +0.68%  '
$ws.Range("E48").Value = '  -0.37%  '
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.48'
$ws.Range("E51").Value = '  -3.46%  '
